$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant_vocab_mapping")

# Change B17 (eex_value for field_collection_field) from numeric 1566
# to the text "ENERGY DATA.INFO"
$ws.Range("B17").Value = "ENERGY DATA.INFO"
